$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $newValue) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.583.36"
Set-TextValue "E2" "  -1.35%  "
Set-TextValue "D3" "1.850.04"
Set-TextValue "E3" "  -1.08%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "314.20"
Set-TextValue "E5" "  -1.32%  "
Set-TextValue "E6" "  +0.02%  "
Set-TextValue "D7" "0.4248"
Set-TextValue "E7" "  -2.44%  "
Set-TextValue "D8" "0.3637"
Set-TextValue "E8" "  -3.16%  "
Set-TextValue "E9" "  +0.38%  "
Set-TextValue "D10" "0.07297"
Set-TextValue "E10" "  -2.38%  "
Set-TextValue "D11" "0.8733"
Set-TextValue "E11" "  -6.83%  "
Set-TextValue "D12" "20.71"
Set-TextValue "E12" "  -2.55%  "
Set-TextValue "D13" "1.892.41"
Set-TextValue "E13" "  -1.77%  "
Set-TextValue "D14" "5.346"
Set-TextValue "E14" "  -1.70%  "
Set-TextValue "D15" "6.524"
Set-TextValue "E15" "  -3.43%  "
Set-TextValue "D16" "0.06921"
Set-TextValue "E16" "  +0.78%  "
Set-TextValue "D17" "1.003"
Set-TextValue "E17" "  -0.05%  "
Set-TextValue "D18" "78.94"
Set-TextValue "E18" "  -3.07%  "
Set-TextValue "D19" "0.000008872"
Set-TextValue "E19" "  -2.00%  "
Set-TextValue "D20" "1.004"
Set-TextValue "E20" "  +0.11%  "
Set-TextValue "D21" "15.40"
Set-TextValue "E21" "  -2.50%  "
Set-TextValue "D22" "27.612.33"
Set-TextValue "E22" "  -1.24%  "
Set-TextValue "D23" "5.009"
Set-TextValue "E23" "  -2.34%  "
Set-TextValue "D24" "10.62"
Set-TextValue "E24" "  -3.93%  "
Set-TextValue "D25" "2.127.37"
Set-TextValue "E25" "  +1.22%  "
Set-TextValue "D26" "1.985"
Set-TextValue "E26" "  -2.74%  "
Set-TextValue "D27" "153.56"
Set-TextValue "E27" "  +0.38%  "
Set-TextValue "D28" "19.02"
Set-TextValue "E28" "  +2.53%  "
Set-TextValue "D29" "121.54"
Set-TextValue "E29" "  +7.10%  "
Set-TextValue "D30" "5.269"
Set-TextValue "E30" "  -5.76%  "
Set-TextValue "E31" "  +12.16%  "
Set-TextValue "D32" "0.08917"
Set-TextValue "E32" "  -1.15%  "
Set-TextValue "D33" "0.7615"
Set-TextValue "E33" "  -6.07%  "
Set-TextValue "D34" "4.573"
Set-TextValue "E34" "  -4.85%  "
Set-TextValue "D35" "2.941"
Set-TextValue "E35" "  -1.13%  "
Set-TextValue "E36" "  -6.88%  "
Set-TextValue "D37" "1.003"
Set-TextValue "E37" "  +0.01%  "
Set-TextValue "D38" "0.05363"
Set-TextValue "E38" "  -2.69%  "
Set-TextValue "D39" "1.092"
Set-TextValue "E39" "  -2.39%  "
Set-TextValue "D40" "0.01943"
Set-TextValue "E40" "  -1.92%  "
Set-TextValue "D41" "2.810"
Set-TextValue "E41" "  -5.37%  "
Set-TextValue "D42" "6.943"
Set-TextValue "E42" "  -0.76%  "
Set-TextValue "D43" "0.5103"
Set-TextValue "E43" "  -3.15%  "
Set-TextValue "D44" "0.1649"
Set-TextValue "E44" "  -2.85%  "
Set-TextValue "D45" "8.290"
Set-TextValue "E45" "  -5.62%  "
Set-TextValue "D46" "0.06560"
Set-TextValue "E46" "  -2.79%  "
Set-TextValue "D49" "104.50"
Set-TextValue "E49" "  -2.22%  "
Set-TextValue "E50" "  +0.04%  "
Set-TextValue "D51" "1.625"
Set-TextValue "E51" "  -2.85%  "

# Rows 47 and 48 swap coin identity (Decentraland <-> EnergySwap) with updated values
Set-TextValue "B47" "Decentraland"
Set-TextValue "C47" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D47" "0.4750"
Set-TextValue "E47" "  -2.86%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "10.40"
Set-TextValue "E48" "  -0.94%  "
